$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 282.63635
$ws.Range("I2").Value = 106.8
$ws.Range("J2").Value = 429.16666
$ws.Range("K2").Value = 106.8
$ws.Range("L2").Value = 429.16666
$ws.Range("M2").Value = 6.200000000000003
$ws.Range("N2").Value = -655.16666
$ws.Range("H9").Value = 214.2
$ws.Range("I9").Value = 262.75
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 262.75
$ws.Range("L9").Value = 20
$ws.Range("M9").Value = -93.75
$ws.Range("N9").Value = -358
$ws.Range("H15").Value = 3105.2097
$ws.Range("I15").Value = 3105.2097
$ws.Range("K15").Value = 9315.6291
$ws.Range("M15").Value = -9146.6291
$ws.Range("H33").Value = 153.83333
$ws.Range("I33").Value = 153.83333
$ws.Range("K33").Value = 153.83333
$ws.Range("M33").Value = 75.16667000000001
$ws.Range("H70").Value = 5267.7144
$ws.Range("I70").Value = 1625
$ws.Range("J70").Value = 7999.75
$ws.Range("K70").Value = 4875
$ws.Range("L70").Value = 23999.25
$ws.Range("M70").Value = -4605
$ws.Range("N70").Value = -24539.25
$ws.Range("H73").Value = 5267.7144
$ws.Range("I73").Value = 1625
$ws.Range("J73").Value = 7999.75
$ws.Range("K73").Value = 4875
$ws.Range("L73").Value = 23999.25
$ws.Range("M73").Value = -3939
$ws.Range("N73").Value = -25871.25
$ws.Range("H86").Value = 1775.5
$ws.Range("I86").Value = 1701
$ws.Range("J86").Value = 1999
$ws.Range("K86").Value = 1701
$ws.Range("L86").Value = 1999
$ws.Range("M86").Value = -578
$ws.Range("N86").Value = -4245
$ws.Range("H88").Value = 1679.6666
$ws.Range("I88").Value = 2127
$ws.Range("J88").Value = 1456
$ws.Range("K88").Value = 2127
$ws.Range("L88").Value = 1456
$ws.Range("M88").Value = -1721
$ws.Range("N88").Value = -2268
$ws.Range("H89").Value = 1775.5
$ws.Range("I89").Value = 1701
$ws.Range("J89").Value = 1999
$ws.Range("K89").Value = 8505
$ws.Range("L89").Value = 9995
$ws.Range("M89").Value = -2889
$ws.Range("N89").Value = -21227
$ws.Range("H91").Value = 1679.6666
$ws.Range("I91").Value = 2127
$ws.Range("J91").Value = 1456
$ws.Range("K91").Value = 2127
$ws.Range("L91").Value = 1456
$ws.Range("M91").Value = -723
$ws.Range("N91").Value = -4264
$ws.Range("H132").Value = 12375.046
$ws.Range("I132").Value = 11571.105
$ws.Range("K132").Value = 34713.315
$ws.Range("M132").Value = -32183.315

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 59899.2
$ws.Range("J114").Value = 59899.2
$ws.Range("L114").Value = 59899.2
$ws.Range("N114").Value = -68577.2
$ws.Range("H132").Value = 4351.6665
$ws.Range("I132").Value = 2574.25
$ws.Range("K132").Value = 7722.75
$ws.Range("M132").Value = -5192.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7090.4546
$ws.Range("I86").Value = 3000.5
$ws.Range("K86").Value = 3000.5
$ws.Range("M86").Value = -1877.5
$ws.Range("H89").Value = 7090.4546
$ws.Range("I89").Value = 3000.5
$ws.Range("K89").Value = 15002.5
$ws.Range("M89").Value = -9386.5
$ws.Range("H94").Value = 1291.6
$ws.Range("I94").Value = 1291.6
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1291.6
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -840.5999999999999
$ws.Range("N94").ClearContents()
$ws.Range("H134").Value = 2212.5789
$ws.Range("I134").Value = 2212.5789
$ws.Range("K134").Value = 6637.736699999999
$ws.Range("M134").Value = -4102.736699999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 993
$ws.Range("I16").Value = 974.1818
$ws.Range("K16").Value = 974.1818
$ws.Range("M16").Value = -687.1818
$ws.Range("H113").Value = 993
$ws.Range("I113").Value = 974.1818
$ws.Range("K113").Value = 974.1818
$ws.Range("M113").Value = 1195.8182

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4475.483
$ws.Range("I80").Value = 4345.913
$ws.Range("K80").Value = 13037.739
$ws.Range("M80").Value = -12101.739
$ws.Range("H83").Value = 4475.483
$ws.Range("I83").Value = 4345.913
$ws.Range("K83").Value = 39113.217
$ws.Range("M83").Value = -34433.217
$ws.Range("H128").Value = 416665.84
$ws.Range("I128").Value = 416665.84
$ws.Range("K128").Value = 1249997.52
$ws.Range("M128").Value = -1245017.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 4110.5
$ws.Range("I10").Value = 4351.25
$ws.Range("J10").Value = 3950
$ws.Range("K10").Value = 4351.25
$ws.Range("L10").Value = 3950
$ws.Range("M10").Value = -4182.25
$ws.Range("N10").Value = -4288
$ws.Range("H29").Value = 1501.5
$ws.Range("J29").Value = 2253.75
$ws.Range("L29").Value = 2253.75
$ws.Range("N29").Value = -2833.75
$ws.Range("H55").Value = 4850.3
$ws.Range("I55").Value = 5885.8
$ws.Range("K55").Value = 5885.8
$ws.Range("M55").Value = -5558.8
$ws.Range("H70").Value = 10601.4
$ws.Range("I70").Value = 4336.6665
$ws.Range("K70").Value = 4336.6665
$ws.Range("M70").Value = -4066.6665
$ws.Range("H73").Value = 10601.4
$ws.Range("I73").Value = 4336.6665
$ws.Range("K73").Value = 4336.6665
$ws.Range("M73").Value = -3400.6665
$ws.Range("H80").Value = 2662.5
$ws.Range("I80").Value = 2504.5
$ws.Range("J80").Value = 2820.5
$ws.Range("K80").Value = 2504.5
$ws.Range("L80").Value = 2820.5
$ws.Range("M80").Value = -1506.5
$ws.Range("N80").Value = -4816.5
$ws.Range("H83").Value = 2662.5
$ws.Range("I83").Value = 2504.5
$ws.Range("J83").Value = 2820.5
$ws.Range("K83").Value = 12522.5
$ws.Range("L83").Value = 14102.5
$ws.Range("M83").Value = -7530.5
$ws.Range("N83").Value = -24086.5
$ws.Range("H97").Value = 769.5833
$ws.Range("I97").Value = 748.7273
$ws.Range("K97").Value = 748.7273
$ws.Range("M97").Value = -252.7273
$ws.Range("H126").Value = 4500.25
$ws.Range("I126").Value = 3750.5
$ws.Range("K126").Value = 11251.5
$ws.Range("M126").Value = -8781.5
$ws.Range("H132").Value = 66037.82000000001
$ws.Range("I132").Value = 97970.82000000001
$ws.Range("J132").Value = 7494
$ws.Range("K132").Value = 293912.46
$ws.Range("L132").Value = 22482
$ws.Range("M132").Value = -291382.46
$ws.Range("N132").Value = -27542

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 340666.34
$ws.Range("J20").Value = 506499.5
$ws.Range("L20").Value = 506499.5
$ws.Range("N20").Value = -506951.5
$ws.Range("H22").Value = 1999
$ws.Range("I22").Value = 1999
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1999
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1704
$ws.Range("N22").ClearContents()
$ws.Range("H23").Value = 9678.666999999999
$ws.Range("I23").Value = 9678.666999999999
$ws.Range("K23").Value = 9678.666999999999
$ws.Range("M23").Value = -9448.666999999999
$ws.Range("H27").Value = 1999
$ws.Range("I27").Value = 1999
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1999
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1892
$ws.Range("N27").ClearContents()
$ws.Range("H33").Value = 9990
$ws.Range("I33").Value = 9990
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 9990
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9700
$ws.Range("N33").ClearContents()
$ws.Range("H69").Value = 50239.4
$ws.Range("J69").Value = 50239.4
$ws.Range("L69").Value = 50239.4
$ws.Range("N69").Value = -51861.4
$ws.Range("H72").Value = 50239.4
$ws.Range("J72").Value = 50239.4
$ws.Range("L72").Value = 150718.2
$ws.Range("N72").Value = -158830.2
$ws.Range("H93").Value = 1110.125
$ws.Range("J93").Value = 943.3333
$ws.Range("L93").Value = 943.3333
$ws.Range("N93").Value = -3439.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 38050
$ws.Range("J98").Value = 38050
$ws.Range("L98").Value = 38050
$ws.Range("N98").Value = -44040
$ws.Range("H132").Value = 3449.0833
$ws.Range("I132").Value = 2330.125
$ws.Range("J132").Value = 5687
$ws.Range("K132").Value = 6990.375
$ws.Range("L132").Value = 17061
$ws.Range("M132").Value = -4460.375
$ws.Range("N132").Value = -22121
